# Clean up test data in the "Test New Survey Import" sheet:
#  - rename the sheet to "test_new_survey_import 1"
#  - replace the sample drug-name rows with generic FreeText test-question rows
#  - clear out the (now unused) third data row
#  - leave the active selection on E23, matching the author's last edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "test_new_survey_import 1"

# Write column B ("type") first across both rows so the shared-string table
# is populated in the same order as the authored workbook.
$ws.Range("B2").Value = "FreeText"
$ws.Range("B3").Value = "FreeText"

$ws.Range("A2").Value = "fdfuu42a22321c123a8_test"
$ws.Range("C2").Value = "Test question fdfuu42a22321c123a8_test"
$ws.Range("D2").Value = "Test question fdfuu42a22321c123a8_test"

$ws.Range("A3").Value = "fdfzz42a66321c123a8_test"
$ws.Range("C3").Value = "Test question fdfzz42a66321c123a8_test"
$ws.Range("D3").Value = "Test question fdfzz42a66321c123a8_test"

# Row 4 no longer holds a third sample question - clear it entirely.
$ws.Range("A4:D4").ClearContents()

# Match the author's final selection before saving.
$ws.Range("E23").Select()
